# Convert columns D, F, K on the "datos" sheet (rows 2-24) from boolean
# TRUE/FALSE cells into plain numeric 0/1 cells with the new values shown
# in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datos")

# Row => @(D, F, K) numeric values to write.
$rowValues = @{
    2  = @(0, 0, 1)
    3  = @(0, 0, 0)
    4  = @(0, 0, 0)
    5  = @(0, 0, 1)
    6  = @(0, 0, 0)
    7  = @(0, 0, 0)
    8  = @(0, 0, 0)
    9  = @(0, 0, 1)
    10 = @(0, 0, 1)
    11 = @(0, 0, 1)
    12 = @(0, 0, 1)
    13 = @(0, 0, 1)
    14 = @(0, 0, 1)
    15 = @(0, 0, 1)
    16 = @(0, 0, 1)
    17 = @(0, 0, 0)
    18 = @(0, 0, 1)
    19 = @(0, 0, 1)
    20 = @(0, 0, 1)
    21 = @(0, 0, 1)
    22 = @(1, 0, 0)
    23 = @(0, 0, 1)
    24 = @(0, 0, 1)
}

foreach ($row in $rowValues.Keys) {
    $vals = $rowValues[$row]
    $ws.Cells.Item($row, 4).Value = [double]$vals[0]   # column D
    $ws.Cells.Item($row, 6).Value = [double]$vals[1]   # column F
    $ws.Cells.Item($row, 11).Value = [double]$vals[2]  # column K
}
